$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B5:P5").Value = 0
